$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H32").Value = 2346.3333
$ws.Range("J32").Value = 2890.3333
$ws.Range("L32").Value = 2890.3333
$ws.Range("N32").Value = -3542.3333
$ws.Range("H138").Value = 3729.2188
$ws.Range("I138").Value = 4580.273
$ws.Range("J138").Value = 3283.4285
$ws.Range("K138").Value = 13740.819
$ws.Range("L138").Value = 9850.2855
$ws.Range("M138").Value = -8600.819
$ws.Range("N138").Value = -20130.2855
$ws.Range("H141").Value = 6828.4287
$ws.Range("I141").Value = 9038.8
$ws.Range("K141").Value = 27116.4
$ws.Range("M141").Value = -21936.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 734.8421
$ws.Range("I2").Value = 501.33334
$ws.Range("K2").Value = 501.33334
$ws.Range("M2").Value = -388.33334
$ws.Range("H74").Value = 87948.5
$ws.Range("I74").Value = 203076.4
$ws.Range("K74").Value = 203076.4
$ws.Range("M74").Value = -202202.4
$ws.Range("H77").Value = 87948.5
$ws.Range("I77").Value = 203076.4
$ws.Range("K77").Value = 1015382
$ws.Range("M77").Value = -1011014
$ws.Range("H116").Value = 734.8421
$ws.Range("I116").Value = 501.33334
$ws.Range("K116").Value = 501.33334
$ws.Range("M116").Value = 1792.66666
$ws.Range("H123").Value = 76785.71
$ws.Range("I123").Value = 76785.71
$ws.Range("K123").Value = 76785.71
$ws.Range("M123").Value = -71885.71

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 734.8421
$ws.Range("I3").Value = 501.33334
$ws.Range("K3").Value = 501.33334
$ws.Range("M3").Value = -387.33334
$ws.Range("H11").Value = 2482.125
$ws.Range("I11").Value = 810
$ws.Range("J11").Value = 7498.5
$ws.Range("K11").Value = 810
$ws.Range("L11").Value = 7498.5
$ws.Range("M11").Value = -670
$ws.Range("N11").Value = -7778.5
$ws.Range("H20").Value = 1873.1666
$ws.Range("I20").Value = 1702.2858
$ws.Range("K20").Value = 1702.2858
$ws.Range("M20").Value = -1455.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8968
$ws.Range("I17").Value = 8984.857
$ws.Range("J17").Value = 8850
$ws.Range("K17").Value = 8984.857
$ws.Range("L17").Value = 8850
$ws.Range("M17").Value = -8810.857
$ws.Range("N17").Value = -9198
$ws.Range("H86").Value = 61029.285
$ws.Range("I86").Value = 84091
$ws.Range("K86").Value = 84091
$ws.Range("M86").Value = -82968
$ws.Range("H89").Value = 61029.285
$ws.Range("I89").Value = 84091
$ws.Range("K89").Value = 420455
$ws.Range("M89").Value = -414839
$ws.Range("H105").Value = 1220.9286
$ws.Range("I105").Value = 1161
$ws.Range("K105").Value = 1161
$ws.Range("M105").Value = 586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 19999
$ws.Range("I70").Value = 19998
$ws.Range("K70").Value = 59994
$ws.Range("M70").Value = -59679
$ws.Range("H73").Value = 19999
$ws.Range("I73").Value = 19998
$ws.Range("K73").Value = 59994
$ws.Range("M73").Value = -58902
$ws.Range("H75").Value = 948
$ws.Range("I75").Value = 948
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2844
$ws.Range("L75").ClearContents()
$ws.Range("M75").Value = -1846
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 948
$ws.Range("I78").Value = 948
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8532
$ws.Range("L78").ClearContents()
$ws.Range("M78").Value = -3540
$ws.Range("N78").ClearContents()
$ws.Range("H134").Value = 4534.857
$ws.Range("I134").Value = 4534.857
$ws.Range("K134").Value = 13604.571
$ws.Range("M134").Value = -8534.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3310.818
$ws.Range("I80").Value = 2638.1428
$ws.Range("K80").Value = 2638.1428
$ws.Range("M80").Value = -1640.1428
$ws.Range("H83").Value = 3310.818
$ws.Range("I83").Value = 2638.1428
$ws.Range("K83").Value = 13190.714
$ws.Range("M83").Value = -8198.714
$ws.Range("H122").Value = 1829.5416
$ws.Range("I122").Value = 1541.409
$ws.Range("K122").Value = 4624.227000000001
$ws.Range("M122").Value = -2174.227000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2641.6667
$ws.Range("I13").Value = 2641.6667
$ws.Range("K13").Value = 2641.6667
$ws.Range("M13").Value = -2501.6667
$ws.Range("H46").Value = 13413.9
$ws.Range("I46").Value = 18098.684
$ws.Range("K46").Value = 18098.684
$ws.Range("M46").Value = -17910.684
$ws.Range("H82").Value = 2022.3226
$ws.Range("I82").Value = 1153.7142
$ws.Range("J82").Value = 2275.6667
$ws.Range("K82").Value = 1153.7142
$ws.Range("L82").Value = 2275.6667
$ws.Range("M82").Value = -792.7141999999999
$ws.Range("N82").Value = -2997.6667
$ws.Range("H85").Value = 2022.3226
$ws.Range("I85").Value = 1153.7142
$ws.Range("J85").Value = 2275.6667
$ws.Range("K85").Value = 1153.7142
$ws.Range("L85").Value = 2275.6667
$ws.Range("M85").Value = 94.28580000000011
$ws.Range("N85").Value = -4771.6667
$ws.Range("H93").Value = 2066.075
$ws.Range("I93").Value = 1872.2106
$ws.Range("J93").Value = 2241.476
$ws.Range("K93").Value = 1872.2106
$ws.Range("L93").Value = 2241.476
$ws.Range("M93").Value = -624.2106000000001
$ws.Range("N93").Value = -4737.476000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 20062
$ws.Range("J38").Value = 20062
$ws.Range("L38").Value = 20062
$ws.Range("N38").Value = -21008
$ws.Range("H49").Value = 24039
$ws.Range("I49").Value = 24039
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 24039
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -23809
$ws.Range("N49").ClearContents()
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19480
$ws.Range("H75").Value = 64000
$ws.Range("J75").Value = 64000
$ws.Range("L75").Value = 64000
$ws.Range("N75").Value = -65872
$ws.Range("H78").Value = 64000
$ws.Range("J78").Value = 64000
$ws.Range("L78").Value = 192000
$ws.Range("N78").Value = -201360
$ws.Range("H141").Value = 86249.25
$ws.Range("J141").Value = 86249.25
$ws.Range("L141").Value = 86249.25
$ws.Range("N141").Value = -96609.25

Write-Output "Applied all changes"